$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Alice's Data & Plots"
$ws.Range("B10").Value = "Apr. 18"
$ws.Range("C10").Value = 3

$ws.Range("A11").Value = "Meeting (CV)"
$ws.Range("B11").Value = "Apr. 23"
$ws.Range("C11").Value = 1

$ws.Range("F11").Select()
